$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = 0.1846339726339726
$ws.Range("D4").Value = 0.9860784980784981
$ws.Range("E4").Value = 0.9857869337869339
$ws.Range("F4").Value = -0.9996921156921158

$ws.Range("C5").Value = 0.1299310179310179
$ws.Range("D5").Value = 0.07893840693840694
$ws.Range("E5").Value = 0.08092000492000492
$ws.Range("F5").Value = 0.02380664380664381

$ws.Range("C6").Value = -0.03253845253845254
$ws.Range("D6").Value = 0.1981788301788302
$ws.Range("E6").Value = 0.1988854508854509
$ws.Range("F6").Value = -0.07052771852771852

$ws.Range("C7").Value = 0.6748253548253548
$ws.Range("D7").Value = -0.01705280905280905
$ws.Range("E7").Value = -0.01695649695649696
$ws.Range("F7").Value = 0.01700845700845701

$ws.Range("C8").Value = 0.2333666333666334
$ws.Range("D8").Value = -0.01906589506589507
$ws.Range("E8").Value = -0.01902316302316302
$ws.Range("F8").Value = 0.02685227085227086

$ws.Range("C9").Value = -0.005874533874533875
$ws.Range("D9").Value = -0.001835209835209836
$ws.Range("E9").Value = -0.002003198003198003
$ws.Range("F9").Value = 0.002053406053406054

$ws.Range("C10").Value = 0.2312251472251472
$ws.Range("D10").Value = -0.05967889167889168
$ws.Range("E10").Value = -0.05942503142503143
$ws.Range("F10").Value = 0.06388628788628789

$ws.Range("C11").Value = -0.5666374826374827
$ws.Range("D11").Value = -0.01414716214716215
$ws.Range("E11").Value = -0.01417152217152217
$ws.Range("F11").Value = 0.01435358635358636
